# Insert a new data row at row 209 (pushing the existing rows 209-227 down
# to 210-228), and populate it with the new Berenjena price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 209:227 down by one row, creating a blank row 209.
$ws.Rows(209).Insert()

# Populate the new row 209 with the reported values.
$ws.Cells.Item(209, 1).Value  = 9
$ws.Cells.Item(209, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(209, 3).Value  = "Metropolitana"
$ws.Cells.Item(209, 4).Value  = 44578
$ws.Cells.Item(209, 5).Value  = 13
$ws.Cells.Item(209, 6).Value  = 100112001
$ws.Cells.Item(209, 7).Value  = "Berenjena"
$ws.Cells.Item(209, 8).Value  = "Sin especificar"
$ws.Cells.Item(209, 9).Value  = "Primera"
$ws.Cells.Item(209, 10).Value = 160
$ws.Cells.Item(209, 11).Value = 7000
$ws.Cells.Item(209, 12).Value = 8000
$ws.Cells.Item(209, 13).Value = 7500
$ws.Cells.Item(209, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(209, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(209, 16).Value = 150
$ws.Cells.Item(209, 17).Value = 50
$ws.Cells.Item(209, 18).Value = "Hortaliza"
